# Release re-build: the "built on" timestamp embedded in the version
# string moves from "February 03 2026 17.29.55 EST" to
# "February 03 2026 18.05.36 EST". The same substring shows up in a
# handful of cells across both worksheets - the version banner / citation
# on the "About" sheet and the per-row "build_version" column on the
# "Boundaries and methane sources" sheet - so just patch every cell that
# contains the old stamp.

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# "About" sheet: the version banner (A2) and the recommended citation (A6)
foreach ($addr in @("A2", "A6")) {
    $cell = $aboutSheet.Range($addr)
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [string] -and $val.Contains($oldStamp)) {
        $cell.Value = $val.Replace($oldStamp, $newStamp)
    }
}

# "Boundaries and methane sources" sheet: build_version column (S), one
# per feature row (rows 2-8).
foreach ($row in 2..8) {
    $cell = $dataSheet.Range("S$row")
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [string] -and $val.Contains($oldStamp)) {
        $cell.Value = $val.Replace($oldStamp, $newStamp)
    }
}
